$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete the old row 3 (Mediación / 0420194406696, which had the hyperlink on C3).
# This removes that record entirely; row 4 (0420194406718) shifts up to become row 3.
$ws.Rows.Item(3).Delete()

# Step 2: insert a new blank row at position 2. This pushes the former row 2 (ID=1, Juicio)
# down to row 3, and the former row 3 (now holding 0420194406718) down to row 4.
$ws.Rows.Item(2).Insert()

# Step 3: populate the newly inserted row 2 with the new "orden de pago" record.
$ws.Cells.Item(2, 2).Value = "preproducciongestion.segurossura.com.ar"
$ws.Cells.Item(2, 3).Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Cells.Item(2, 4).Value = "lzambonini"
$ws.Cells.Item(2, 5).Value = "silverarrow"

# Match the plain (quote-prefixed text, no fill/border) style used elsewhere in column F/G,
# the same style already present on F4/G4, then write the order number as text so the
# leading zero survives.
$ws.Cells.Item(2, 6).Style = $ws.Cells.Item(4, 6).Style
$ws.Cells.Item(2, 7).Style = $ws.Cells.Item(4, 7).Style
$ws.Cells.Item(2, 6).Value = "'0420172008486"
$ws.Cells.Item(2, 7).Value = "Mediación"

# The row-3 cell that now holds the Juicio record (shifted down from the old row 2)
# must not carry a hyperlink. Remove only that specific hyperlink, leaving the one on
# C4 (the 0420194406718 record) intact.
$target = $ws.Range("C3").Address()
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq $target) {
        $hl.Delete()
    }
}

# Update the active selection to D5, as recorded in the saved view state.
$ws.Range("D5").Select()

$wb.Save()
